$d = $word.ActiveDocument

# The last paragraph in the document body currently reads
# "AggiornamentoManutenzioneNonOrdinariaDittaEsterna" and (unlike the other
# items in this list) has no cached paragraph-mark formatting (no <w:pPr>),
# because no Enter has ever been pressed after it. We are appending one
# more item ("RitardoConsistente") to the list, which means:
#   1. The existing last paragraph stops being the last paragraph, so its
#      paragraph mark now carries the list's usual run formatting
#      (CMU Serif Roman, 10pt) -> a <w:pPr><w:rPr> appears on it.
#   2. A brand new paragraph is appended, holding the new text, formatted
#      the same way as the rest of the list.

$lastPara = $d.Paragraphs.Last
$markRange = $lastPara.Range
$markRange.Collapse(0)            # wdCollapseEnd -> the paragraph mark

# Split off a new (still empty) paragraph after the current last one.
$markRange.InsertParagraphAfter()

# Stamp the paragraph-mark formatting onto the paragraph that used to be
# last (it is now the second-to-last paragraph).
$firstPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$firstMarkRange = $firstPara.Range
$firstMarkRange.Collapse(0)
$firstMarkRange.Font.Name = "CMU Serif Roman"
$firstMarkRange.Font.Size = 10
$firstMarkRange.Font.SizeBi = 10

# Fill in the text of the freshly created last paragraph.
$newLast = $d.Paragraphs.Last
$newRange = $newLast.Range
$newRange.Collapse(0)
$newRange.InsertAfter("RitardoConsistente")
